$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row at 49 (shifts everything from 49..133 down to 50..134) ---
$ws.Rows("49:49").Insert()

# The worksheet table (Table1) does not auto-grow with Rows.Insert(), so
# resize it explicitly to cover the new last row (134).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K134"))

# Restore the calculated-column formula on the newly appended table row
# (134) - Resize() does not copy it down automatically.
$ws.Range("G134").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Pull formatting (number formats / borders / styles) for the brand new
# row 49 from row 45, which already carries the right "data row" style.
$ws.Range("A45:K45").Copy()
$ws.Range("A49:K49").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G49").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Populate the new row 49 (blank-date "UT(0-1-0)" annotation row).
$ws.Range("B49").Value2 = "UT(0-1-0)"
$ws.Range("D49").Value2 = 0.125

# --- Fill in the two new undertime annotations (rows 45 and 47) ---
# Set B47 before B45 so the new shared-string entries are appended in the
# same order as the target workbook (UT(0-1-4) then UT(0-0-32)).
$ws.Range("B47").Value2 = "UT(0-1-4)"
$ws.Range("D47").Value2 = 0.133
$ws.Range("B45").Value2 = "UT(0-0-32)"
$ws.Range("D45").Value2 = 0.067

# --- CONVERTION sheet: switch the little minutes->days calculator from
# 1 day to 32 minutes (used to derive the 0.067 value above) ---
$convert = $wb.Worksheets.Item("CONVERTION")
$convert.Range("E3").ClearContents()
$convert.Range("F3").Value2 = 32

# --- Restore the last active selection on Sheet1 ---
$ws.Range("F45").Select()
